$d = $word.ActiveDocument

# --- Locate the paragraph that currently holds the "_GoBack" bookmark      ---
# --- (the last screenshot before the "Pose Graph Optimization:" heading). ---
$picPara = $null
if ($d.Bookmarks.Exists("_GoBack")) {
    $picPara = $d.Bookmarks.Item("_GoBack").Range.Paragraphs.First
    # Remove the old bookmark from the picture paragraph; it is re-created
    # on the new paragraph inserted below.
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Insert a brand new paragraph right after the picture paragraph.
$newParaRange = $picPara.Range.InsertParagraphAfter()

# Re-fetch the freshly created (still empty) paragraph and fill it in with
# three separate runs plus the "_GoBack" bookmark, via a raw WordprocessingML
# fragment so the run boundaries come out exactly as intended.
$insertedPara = $d.Paragraphs.Item($picPara.Index + 1)
$insertedRange = $insertedPara.Range

$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:r><w:t>Question for</w:t></w:r>' + `
              '<w:r><w:t xml:space="preserve"> bundle adjustment: how </w:t></w:r>' + `
              '<w:r><w:t>do we obtain an initial guess?</w:t></w:r>' + `
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
              '<w:bookmarkEnd w:id="0"/>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

[void]$insertedRange.InsertXML($xmlFragment)

# --- Drop the stray <w:lastRenderedPageBreak/> immediately before "LSD SLAM" ---
[void]$d.Content.Find.Execute("LSD SLAM", $true, $false, $false, $false, $false,
                               $true, 1, $false, "LSD SLAM", 2)

Write-Output "done"
